$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet by copying "2022-Q2" and placing
#    the copy immediately before it (so tab order becomes
#    总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3, 2021-Q2, 2021-Q1, 2020-Q4).
# ------------------------------------------------------------------
$srcQ2 = $wb.Worksheets.Item("2022-Q2")
$srcQ2.Copy($srcQ2, $null)
$newQ3 = $wb.Worksheets.Item(2)
$newQ3.Name = "2022-Q3"

# ------------------------------------------------------------------
# 2) Overwrite the new sheet's data rows with the 2022-Q3 figures.
#    Columns D/E/F/G are stored as text (same as the rest of the
#    workbook), H is a real number.
# ------------------------------------------------------------------
$newQ3.Range("D2").Value = "'1.58"
$newQ3.Range("D2").Style = "Normal"
$newQ3.Range("E2").Value = "'78.58"
$newQ3.Range("E2").Style = "Normal"
$newQ3.Range("F2").Value = "'3.74"
$newQ3.Range("F2").Style = "Normal"
$newQ3.Range("G2").Value = "'0.0591"
$newQ3.Range("G2").Style = "Normal"
$newQ3.Range("H2").Value = 8

$newQ3.Range("D3").Value = "'1.58"
$newQ3.Range("D3").Style = "Normal"
$newQ3.Range("E3").Value = "'78.58"
$newQ3.Range("E3").Style = "Normal"
$newQ3.Range("F3").Value = "'3.74"
$newQ3.Range("F3").Style = "Normal"
$newQ3.Range("G3").Value = "'0.0591"
$newQ3.Range("G3").Style = "Normal"
$newQ3.Range("H3").Value = 8

# ------------------------------------------------------------------
# 3) Insert a new top data row in the "总计" (summary) sheet for the
#    2022-Q3 quarter, pushing the existing rows down by one.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.12

# Match the plain-number-column style ("s=2") used by the other rows
# in column A.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
